# Append two new daily-log rows (2020-08-27 and 2020-08-28) to the
# "out_vars" sheet, right after the existing last row (88).
#
# Column A holds dates stored as plain text (e.g. "2020-06-01"), not real
# Excel date serials. To keep new entries consistent with the existing
# column, the cells are pre-formatted as Text ("@") so Excel doesn't
# auto-convert the typed string into a date serial number; the format is
# then reset back to General (matching the rest of the column) once the
# text values are safely in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A89:A90").NumberFormat = "@"

$ws.Range("A89").Value = "2020-08-27"
$ws.Range("B89").Value = 579914
$ws.Range("C89").Value = 643265
$ws.Range("D89").Value = 81597
$ws.Range("E89").Value = 62594
$ws.Range("F89").Value = 25.64

$ws.Range("A90").Value = "2020-08-28"
$ws.Range("B90").Value = 585738
$ws.Range("C90").Value = 650862
$ws.Range("D90").Value = 83357
$ws.Range("E90").Value = 63146
$ws.Range("F90").Value = 25.57

$ws.Range("A89:A90").NumberFormat = "general"

# Mirror the author's final UI state: cursor parked on the first blank row
# below the new data, scrolled down so the newly-entered rows are visible.
$ws.Range("A91:XFD91").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A82")
